$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated crypto price/volume data (and the row-48 BabyDogeCoin insertion
# that shifted Algorand/Cronos/USDD/EnergySwap down) cell by cell.
# NumberFormat is forced to text ("@") before the assignment and restored to the
# sheet default ("Normal" style) right after, so numeric-looking strings such as
# "1.00" or "211.38" are preserved verbatim instead of being coerced into numbers.

$c = $ws.Range('D2')
$c.NumberFormat = "@"
$c.Value = '26.663.62'
$c.Style = "Normal"

$c = $ws.Range('E2')
$c.NumberFormat = "@"
$c.Value = '  +0.05%  '
$c.Style = "Normal"

$c = $ws.Range('D3')
$c.NumberFormat = "@"
$c.Value = '1.596.96'
$c.Style = "Normal"

$c = $ws.Range('E3')
$c.NumberFormat = "@"
$c.Value = '  -0.11%  '
$c.Style = "Normal"

$c = $ws.Range('E4')
$c.NumberFormat = "@"
$c.Value = '  +0.16%  '
$c.Style = "Normal"

$c = $ws.Range('D5')
$c.NumberFormat = "@"
$c.Value = '211.38'
$c.Style = "Normal"

$c = $ws.Range('E5')
$c.NumberFormat = "@"
$c.Value = '  +0.27%  '
$c.Style = "Normal"

$c = $ws.Range('D6')
$c.NumberFormat = "@"
$c.Value = '0.512'
$c.Style = "Normal"

$c = $ws.Range('E6')
$c.NumberFormat = "@"
$c.Value = '  +0.00%  '
$c.Style = "Normal"

$c = $ws.Range('E7')
$c.NumberFormat = "@"
$c.Value = '  +0.11%  '
$c.Style = "Normal"

$c = $ws.Range('E8')
$c.NumberFormat = "@"
$c.Value = '  +0.04%  '
$c.Style = "Normal"

$c = $ws.Range('E9')
$c.NumberFormat = "@"
$c.Value = '  +0.57%  '
$c.Style = "Normal"

$c = $ws.Range('D10')
$c.NumberFormat = "@"
$c.Value = '19.48'
$c.Style = "Normal"

$c = $ws.Range('E10')
$c.NumberFormat = "@"
$c.Value = '  -0.96%  '
$c.Style = "Normal"

$c = $ws.Range('D11')
$c.NumberFormat = "@"
$c.Value = '0.0841'
$c.Style = "Normal"

$c = $ws.Range('E11')
$c.NumberFormat = "@"
$c.Value = '  +0.37%  '
$c.Style = "Normal"

$c = $ws.Range('D12')
$c.NumberFormat = "@"
$c.Value = '1.821.48'
$c.Style = "Normal"

$c = $ws.Range('E12')
$c.NumberFormat = "@"
$c.Value = '  -0.07%  '
$c.Style = "Normal"

$c = $ws.Range('D13')
$c.NumberFormat = "@"
$c.Value = '1.583.33'
$c.Style = "Normal"

$c = $ws.Range('E13')
$c.NumberFormat = "@"
$c.Value = '  -1.33%  '
$c.Style = "Normal"

$c = $ws.Range('E14')
$c.NumberFormat = "@"
$c.Value = '  +0.15%  '
$c.Style = "Normal"

$c = $ws.Range('D15')
$c.NumberFormat = "@"
$c.Value = '0.523'
$c.Style = "Normal"

$c = $ws.Range('E15')
$c.NumberFormat = "@"
$c.Value = '  +0.36%  '
$c.Style = "Normal"

$c = $ws.Range('D16')
$c.NumberFormat = "@"
$c.Value = '65.05'
$c.Style = "Normal"

$c = $ws.Range('E16')
$c.NumberFormat = "@"
$c.Value = '  +0.28%  '
$c.Style = "Normal"

$c = $ws.Range('D17')
$c.NumberFormat = "@"
$c.Value = '26.640.82'
$c.Style = "Normal"

$c = $ws.Range('E17')
$c.NumberFormat = "@"
$c.Value = '  +0.01%  '
$c.Style = "Normal"

$c = $ws.Range('D18')
$c.NumberFormat = "@"
$c.Value = '0.0₃0738'
$c.Style = "Normal"

$c = $ws.Range('E18')
$c.NumberFormat = "@"
$c.Value = '  +1.17%  '
$c.Style = "Normal"

$c = $ws.Range('E19')
$c.NumberFormat = "@"
$c.Value = '  +0.20%  '
$c.Style = "Normal"

$c = $ws.Range('D20')
$c.NumberFormat = "@"
$c.Value = '209.02'
$c.Style = "Normal"

$c = $ws.Range('E20')
$c.NumberFormat = "@"
$c.Value = '  -0.10%  '
$c.Style = "Normal"

$c = $ws.Range('D21')
$c.NumberFormat = "@"
$c.Value = '7.06'
$c.Style = "Normal"

$c = $ws.Range('E21')
$c.NumberFormat = "@"
$c.Value = '  +4.42%  '
$c.Style = "Normal"

$c = $ws.Range('E22')
$c.NumberFormat = "@"
$c.Value = '  +0.40%  '
$c.Style = "Normal"

$c = $ws.Range('E23')
$c.NumberFormat = "@"
$c.Value = '  +2.64%  '
$c.Style = "Normal"

$c = $ws.Range('D24')
$c.NumberFormat = "@"
$c.Value = '9.00'
$c.Style = "Normal"

$c = $ws.Range('E24')
$c.NumberFormat = "@"
$c.Value = '  +1.08%  '
$c.Style = "Normal"

$c = $ws.Range('D25')
$c.NumberFormat = "@"
$c.Value = '143.82'
$c.Style = "Normal"

$c = $ws.Range('E25')
$c.NumberFormat = "@"
$c.Value = '  -1.41%  '
$c.Style = "Normal"

$c = $ws.Range('E26')
$c.NumberFormat = "@"
$c.Value = '  +0.15%  '
$c.Style = "Normal"

$c = $ws.Range('E27')
$c.NumberFormat = "@"
$c.Value = '  -1.79%  '
$c.Style = "Normal"

$c = $ws.Range('E28')
$c.NumberFormat = "@"
$c.Value = '  -1.13%  '
$c.Style = "Normal"

$c = $ws.Range('D29')
$c.NumberFormat = "@"
$c.Value = '15.31'
$c.Style = "Normal"

$c = $ws.Range('E29')
$c.NumberFormat = "@"
$c.Value = '  +0.15%  '
$c.Style = "Normal"

$c = $ws.Range('D30')
$c.NumberFormat = "@"
$c.Value = '0.0514'
$c.Style = "Normal"

$c = $ws.Range('E30')
$c.NumberFormat = "@"
$c.Value = '  +1.69%  '
$c.Style = "Normal"

$c = $ws.Range('E31')
$c.NumberFormat = "@"
$c.Value = '  -0.03%  '
$c.Style = "Normal"

$c = $ws.Range('E32')
$c.NumberFormat = "@"
$c.Value = '  -0.10%  '
$c.Style = "Normal"

$c = $ws.Range('E33')
$c.NumberFormat = "@"
$c.Value = '  +0.68%  '
$c.Style = "Normal"

$c = $ws.Range('D34')
$c.NumberFormat = "@"
$c.Value = '1.291.19'
$c.Style = "Normal"

$c = $ws.Range('E34')
$c.NumberFormat = "@"
$c.Value = '  -0.30%  '
$c.Style = "Normal"

$c = $ws.Range('D35')
$c.NumberFormat = "@"
$c.Value = '0.619'
$c.Style = "Normal"

$c = $ws.Range('E35')
$c.NumberFormat = "@"
$c.Value = '  -6.32%  '
$c.Style = "Normal"

$c = $ws.Range('E36')
$c.NumberFormat = "@"
$c.Value = '  +0.57%  '
$c.Style = "Normal"

$c = $ws.Range('E37')
$c.NumberFormat = "@"
$c.Value = '  -0.03%  '
$c.Style = "Normal"

$c = $ws.Range('E38')
$c.NumberFormat = "@"
$c.Value = '  -0.30%  '
$c.Style = "Normal"

$c = $ws.Range('E39')
$c.NumberFormat = "@"
$c.Value = '  -1.34%  '
$c.Style = "Normal"

$c = $ws.Range('E40')
$c.NumberFormat = "@"
$c.Value = '  +15.02%  '
$c.Style = "Normal"

$c = $ws.Range('E41')
$c.NumberFormat = "@"
$c.Value = '  +1.28%  '
$c.Style = "Normal"

$c = $ws.Range('E42')
$c.NumberFormat = "@"
$c.Value = '  -0.55%  '
$c.Style = "Normal"

$c = $ws.Range('D43')
$c.NumberFormat = "@"
$c.Value = '0.783'
$c.Style = "Normal"

$c = $ws.Range('E43')
$c.NumberFormat = "@"
$c.Value = '  -0.49%  '
$c.Style = "Normal"

$c = $ws.Range('D44')
$c.NumberFormat = "@"
$c.Value = '63.23'
$c.Style = "Normal"

$c = $ws.Range('E44')
$c.NumberFormat = "@"
$c.Value = '  -0.90%  '
$c.Style = "Normal"

$c = $ws.Range('D45')
$c.NumberFormat = "@"
$c.Value = '1.732.67'
$c.Style = "Normal"

$c = $ws.Range('E45')
$c.NumberFormat = "@"
$c.Value = '  -0.18%  '
$c.Style = "Normal"

$c = $ws.Range('D46')
$c.NumberFormat = "@"
$c.Value = '91.10'
$c.Style = "Normal"

$c = $ws.Range('E46')
$c.NumberFormat = "@"
$c.Value = '  +0.98%  '
$c.Style = "Normal"

$c = $ws.Range('E47')
$c.NumberFormat = "@"
$c.Value = '  -3.07%  '
$c.Style = "Normal"

$c = $ws.Range('B48')
$c.NumberFormat = "@"
$c.Value = 'BabyDogeCoin'
$c.Style = "Normal"

$c = $ws.Range('C48')
$c.NumberFormat = "@"
$c.Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$c.Style = "Normal"

$c = $ws.Range('D48')
$c.NumberFormat = "@"
$c.Value = '0.0₆0104'
$c.Style = "Normal"

$c = $ws.Range('E48')
$c.NumberFormat = "@"
$c.Value = '  -1.40%  '
$c.Style = "Normal"

$c = $ws.Range('B49')
$c.NumberFormat = "@"
$c.Value = 'Algorand'
$c.Style = "Normal"

$c = $ws.Range('C49')
$c.NumberFormat = "@"
$c.Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$c.Style = "Normal"

$c = $ws.Range('D49')
$c.NumberFormat = "@"
$c.Value = '0.101'
$c.Style = "Normal"

$c = $ws.Range('E49')
$c.NumberFormat = "@"
$c.Value = '  +0.79%  '
$c.Style = "Normal"

$c = $ws.Range('B50')
$c.NumberFormat = "@"
$c.Value = 'Cronos'
$c.Style = "Normal"

$c = $ws.Range('C50')
$c.NumberFormat = "@"
$c.Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$c.Style = "Normal"

$c = $ws.Range('D50')
$c.NumberFormat = "@"
$c.Value = '0.0509'
$c.Style = "Normal"

$c = $ws.Range('E50')
$c.NumberFormat = "@"
$c.Value = '  +0.94%  '
$c.Style = "Normal"

$c = $ws.Range('B51')
$c.NumberFormat = "@"
$c.Value = 'USDD'
$c.Style = "Normal"

$c = $ws.Range('C51')
$c.NumberFormat = "@"
$c.Value = 'https://coinranking.com/coin/z2PZIKQL7+usdd-usdd'
$c.Style = "Normal"

$c = $ws.Range('D51')
$c.NumberFormat = "@"
$c.Value = '1.00'
$c.Style = "Normal"

$c = $ws.Range('E51')
$c.NumberFormat = "@"
$c.Value = '  +0.08%  '
$c.Style = "Normal"
